$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table shrinks from 6 data rows to 5 - drop the last one.
$ws.Rows.Item(6).Delete() | Out-Null

# --- Row 1 (headers) stay the same content, just re-assert them ---
$ws.Range("A1").Value = "CasosDePrueba"
$ws.Range("B1").Value = "dato 1"
$ws.Range("C1").Value = "dato 2"
$ws.Range("D1").Value = "dato 3"

# --- Row 2 ---
$ws.Range("A2").Value = "CP001_loginInvalidEmail"
$ws.Range("B2").Value = "qweqweew"
$ws.Range("C2").Value = "ee51165"
$ws.Range("D2").Value = "Invalid email address."

# --- Row 3 (new test case: invalid login e-mail) ---
$ws.Range("B3").Value = "fakemail@gmail.com"
$ws.Range("D3").Value = "Authentication failed."
$ws.Range("C3").Value = "asdffgr2"
$ws.Range("A3").Value = "CP002_loginEmail"
$ws.Range("A3").Font.Underline = $true

# --- Row 4 (new test case: successful login) ---
$ws.Range("B4").Value = "mailtestautomation001@gmail.com"
$ws.Range("C4").Value = "admin123"
$ws.Range("D4").Value = "Welcome to your account. Here you can manage all of your personal information and orders."
$ws.Range("A4").Value = "CP003_logInSuccessfully"

# --- Row 5 (new test case: add product to wishlist) ---
$ws.Range("B5").Value = "mailtestautomation001@gmail.com"
$ws.Range("C5").Value = "admin123"
$ws.Range("A5").Value = "CP004_AddProductToWishlist"
$ws.Range("D5").Value = "Added to your wishlist."
$ws.Range("D5").Font.Underline = $false

# --- Hyperlinks on the e-mail cells (B3:B5) ---
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:fakemail@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:mailtestautomation001@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:mailtestautomation001@gmail.com") | Out-Null

# --- Underline the two "admin123" credential cells, with no border ---
$ws.Range("C4").Font.Underline = $true
$ws.Range("C4").Borders.LineStyle = -4142
$ws.Range("C5").Font.Underline = $true
$ws.Range("C5").Borders.LineStyle = -4142

# --- Column D is now wide enough to fit the long messages ---
$ws.Columns.Item(4).ColumnWidth = 83.3

# --- Selection moves to D5 to mirror the saved state ---
$ws.Range("D5").Select() | Out-Null
